# arvo_id 40096184 (row 5) is corrected to reflect that the fix actually
# cleared the crash.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# first_patch_removed_crash: FALSE -> TRUE
$ws.Range("B5").Value = $true

# notes: drop the stale "Cost spiked in time and tokens." remark now that
# the run is understood to have cleared the crash.
$ws.Range("N5").Value = "Patched different file than GT, but crash still cleared."

# Reflect the reviewer's on-screen scroll/selection state at save time
# (the view had scrolled right so column F was leftmost, with N6 selected).
$excel.Goto($ws.Range("N6"), $true)
$excel.ActiveWindow.ScrollColumn = 6
